$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.430.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("D3").Value = "'3.175.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.98%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'587.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.34%  "
$ws.Range("D6").Value = "'135.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.20%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'3.173.72"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.91%  "
$ws.Range("E9").Value = "  -2.79%  "
$ws.Range("D10").Value = "'0.141"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.06%  "
$ws.Range("D11").Value = "'5.25"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.33%  "
$ws.Range("E12").Value = "  -3.94%  "
$ws.Range("E13").Value = "  -5.40%  "
$ws.Range("D14").Value = "'33.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.24%  "
$ws.Range("D15").Value = "'3.697.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.03%  "
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("D17").Value = "'3.174.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.05%  "
$ws.Range("D18").Value = "'62.437.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.33%  "
$ws.Range("E19").Value = "  -5.10%  "
$ws.Range("D20").Value = "'455.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.38%  "
$ws.Range("D21").Value = "'13.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.82%  "
$ws.Range("E22").Value = "  -4.16%  "
$ws.Range("E23").Value = "  -4.78%  "
$ws.Range("D24").Value = "'13.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("D25").Value = "'83.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.80%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -3.09%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").Value = "'6.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.49%  "
$ws.Range("D30").Value = "'7.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.02%  "
$ws.Range("E31").Value = "  -6.86%  "
$ws.Range("D32").Value = "'27.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.96%  "
$ws.Range("D33").Value = "'0.104"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.58%  "
$ws.Range("E34").Value = "  -6.15%  "
$ws.Range("E35").Value = "  -5.92%  "
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("D37").Value = "'51.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.36%  "
$ws.Range("E38").Value = "  -7.30%  "
$ws.Range("D39").Value = "'0.0384"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.41%  "
$ws.Range("D40").Value = "'410.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.42%  "
$ws.Range("D41").Value = "'2.67"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.41%  "
$ws.Range("D42").Value = "'2.871.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.57%  "
$ws.Range("E43").Value = "  -4.72%  "
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("D45").Value = "'36.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.97%  "
$ws.Range("D46").Value = "'0.248"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.27%  "
$ws.Range("D49").Value = "'124.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("D50").Value = "'25.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.93%  "
$ws.Range("E51").Value = "  -3.61%  "

# Row 47/48: USDe and Fetch.AI swap ranking positions
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "'0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").Value = "'2.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.88%  "
